{"js": "// The document contains a single table of simple arithmetic problems\n// (\"50-31=\", \"17+54=\", ...). The edit replaces the text of every\n// problem cell (in reading order: row by row, left to right) with a\n// new expression, while leaving the table's structure (rows/columns),\n// formatting and the date paragraph above it untouched.\nconst newValues = [\n  [\"46+17=\", \"12-3=\", \"89+8=\", \"48+44=\", \"12-5=\"],\n  [\"16+75=\", \"29+8=\", \"9+6=\", \"85-69=\", \"37+19=\"],\n  [\"32+59=\", \"73-18=\", \"28+38=\", \"7+67=\", \"70-19=\"],\n  [\"6+66=\", \"25+36=\", \"71-52=\", \"79+18=\", \"44+9=\"],\n  [\"28+15=\", \"32-5=\", \"28+13=\", \"80-79=\", \"59+39=\"],\n  [\"31-25=\", \"90-55=\", \"55+19=\", \"62-25=\", \"82-54=\"],\n  [\"28+15=\", \"92-16=\", \"8+57=\", \"32-14=\", \"7+39=\"],\n  [\"23-6=\", \"47-19=\", \"92-83=\", \"9+65=\", \"86-9=\"],\n  [\"70-34=\", \"64-27=\", \"29+6=\", \"8+73=\", \"41-22=\"],\n  [\"20-14=\", \"79+13=\", \"84-26=\", \"61-26=\", \"72-53=\"],\n  [\"7+66=\", \"6+85=\", \"33-28=\", \"12+79=\", \"16+35=\"],\n  [\"26+27=\", \"19+17=\", \"83-6=\", \"84-55=\", \"91-83=\"],\n  [\"17+77=\", \"76-68=\", \"23-4=\", \"50-9=\", \"16+38=\"],\n  [\"4+28=\", \"36+6=\", \"72-36=\", \"67-29=\", \"81-65=\"],\n  [\"32-17=\", \"35+39=\", \"58-9=\", \"76+6=\", \"18+78=\"],\n  [\"67+5=\", \"69+14=\", \"7+34=\", \"34+48=\", \"80-73=\"],\n  [\"55-26=\", \"9+19=\", \"64-6=\", \"23+48=\", \"3+48=\"],\n  [\"58+37=\", \"55-18=\", \"52+19=\", \"61-19=\", \"90-45=\"],\n  [\"40-31=\", \"15+47=\", \"33+58=\", \"39+8=\", \"77+14=\"],\n  [\"85-67=\", \"82-33=\", \"90-34=\", \"18+27=\", \"8+83=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\n// Write row by row so a shape mismatch fails loudly instead of silently\n// mis-mapping cells.\nif (table.rowCount !== newValues.length) {\n  throw new Error(\n    `Expected ${newValues.length} rows, found ${table.rowCount}`\n  );\n}\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# The document contains a single table of simple arithmetic problems\n# (\"50-31=\", \"17+54=\", ...). The edit replaces the text of every\n# problem cell (in reading order: row by row, left to right) with a\n# new expression, while leaving the table's structure (rows/columns),\n# formatting and the date paragraph above it untouched.\n$newValues = @(\n  @(\"46+17=\", \"12-3=\", \"89+8=\", \"48+44=\", \"12-5=\"),\n  @(\"16+75=\", \"29+8=\", \"9+6=\", \"85-69=\", \"37+19=\"),\n  @(\"32+59=\", \"73-18=\", \"28+38=\", \"7+67=\", \"70-19=\"),\n  @(\"6+66=\", \"25+36=\", \"71-52=\", \"79+18=\", \"44+9=\"),\n  @(\"28+15=\", \"32-5=\", \"28+13=\", \"80-79=\", \"59+39=\"),\n  @(\"31-25=\", \"90-55=\", \"55+19=\", \"62-25=\", \"82-54=\"),\n  @(\"28+15=\", \"92-16=\", \"8+57=\", \"32-14=\", \"7+39=\"),\n  @(\"23-6=\", \"47-19=\", \"92-83=\", \"9+65=\", \"86-9=\"),\n  @(\"70-34=\", \"64-27=\", \"29+6=\", \"8+73=\", \"41-22=\"),\n  @(\"20-14=\", \"79+13=\", \"84-26=\", \"61-26=\", \"72-53=\"),\n  @(\"7+66=\", \"6+85=\", \"33-28=\", \"12+79=\", \"16+35=\"),\n  @(\"26+27=\", \"19+17=\", \"83-6=\", \"84-55=\", \"91-83=\"),\n  @(\"17+77=\", \"76-68=\", \"23-4=\", \"50-9=\", \"16+38=\"),\n  @(\"4+28=\", \"36+6=\", \"72-36=\", \"67-29=\", \"81-65=\"),\n  @(\"32-17=\", \"35+39=\", \"58-9=\", \"76+6=\", \"18+78=\"),\n  @(\"67+5=\", \"69+14=\", \"7+34=\", \"34+48=\", \"80-73=\"),\n  @(\"55-26=\", \"9+19=\", \"64-6=\", \"23+48=\", \"3+48=\"),\n  @(\"58+37=\", \"55-18=\", \"52+19=\", \"61-19=\", \"90-45=\"),\n  @(\"40-31=\", \"15+47=\", \"33+58=\", \"39+8=\", \"77+14=\"),\n  @(\"85-67=\", \"82-33=\", \"90-34=\", \"18+27=\", \"8+83=\")\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nif ($t.Rows.Count -ne $newValues.Length) {\n  throw \"Expected $($newValues.Length) rows, found $($t.Rows.Count)\"\n}\n\nfor ($r = 0; $r -lt $newValues.Length; $r++) {\n  $row = $newValues[$r]\n  for ($c = 0; $c -lt $row.Length; $c++) {\n    $t.Cell($r + 1, $c + 1).Range.Text = $row[$c]\n  }\n}\n"}
